$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = 45209
$ws.Range("A62").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B62").Value = "21:15"
$ws.Range("C62").Value = 2563.32
$ws.Range("D62").Value = "amazon"
$ws.Range("E62").Value = "preto"

$ws.Range("A63").Value = 45209
$ws.Range("A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B63").Value = "21:16"
$ws.Range("C63").Value = 2563
$ws.Range("D63").Value = "mercado livre"
$ws.Range("E63").Value = "preto"
